$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-02-18"

# Update the header label in I1 (shared string) to match the new date
$ws.Range("I1").Value = "2022 (through 02-18)"

# Update the data values for the new day's data (2022-02-26 commit, February & Total rows)
$ws.Range("I3").Value = 85
$ws.Range("I14").Value = 245
